{"js": "// Fixed #198 Maximum size of 255 char for fields may be a problem for the\n// list of nsURIs to use.\n//\n// The bookmark (\"_GoBack\") that used to sit at the end of the first\n// paragraph is moved to its own (previously empty) trailing paragraph,\n// and a run of the field instruction text that was split across several\n// <w:instrText> runs gets merged into a single run.\n\nconst paras = context.document.body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\n\nconst firstPara = paras.items[0];\nconst fieldPara = paras.items[1];\nconst lastPara = paras.items[paras.items.length - 1];\n\nconst nsWord = 'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"';\n\n// 1) Paragraph 1: \"A simple demonstration of a query :\" -- drop the\n//    bookmark that used to live here (it is recreated below, in the\n//    trailing paragraph).\nconst firstParaOoxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document ' + nsWord + '><w:body>' +\n  '<w:p>' +\n  '<w:r><w:t xml:space=\"preserve\">A simple </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:t>demonstration</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> of a </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:t>query</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t>\\u00a0:</w:t></w:r>' +\n  '</w:p>' +\n  '</w:body></w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\nfirstPara.insertOoxml(firstParaOoxml, Word.InsertLocation.replace);\n\n// 2) Paragraph 2: the field-code paragraph. Content/text is unchanged;\n//    only the run split of the instruction text changes: the runs\n//    holding \",\", \"''.sampleTable()\", \",\", \"''.sampleTable()\", \"}\" are\n//    merged into a single run \",''.sampleTable(),''.sampleTable()}\".\nconst fieldParaOoxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document ' + nsWord + '><w:body>' +\n  '<w:p>' +\n  '<w:r><w:fldChar w:fldCharType=\"begin\"/></w:r>' +\n  '<w:r><w:instrText xml:space=\"preserve\"> </w:instrText></w:r>' +\n  '<w:r><w:instrText>m</w:instrText></w:r>' +\n  '<w:r><w:instrText>:</w:instrText></w:r>' +\n  '<w:r><w:instrText>Sequence</w:instrText></w:r>' +\n  '<w:r><w:instrText>{</w:instrText></w:r>' +\n  \"<w:r><w:instrText>''.</w:instrText></w:r>\" +\n  '<w:r><w:instrText>sample</w:instrText></w:r>' +\n  '<w:r><w:instrText>Table()</w:instrText></w:r>' +\n  \"<w:r><w:instrText>,''.sampleTable(),''.sampleTable()}</w:instrText></w:r>\" +\n  '<w:r><w:instrText xml:space=\"preserve\"> </w:instrText></w:r>' +\n  '<w:r><w:fldChar w:fldCharType=\"end\"/></w:r>' +\n  '</w:p>' +\n  '</w:body></w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\nfieldPara.insertOoxml(fieldParaOoxml, Word.InsertLocation.replace);\n\n// 3) The final (previously empty) paragraph now hosts the \"_GoBack\"\n//    bookmark that used to be attached to the first paragraph.\nconst bookmarkParaOoxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document ' + nsWord + '><w:body>' +\n  '<w:p>' +\n  '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n  '<w:bookmarkEnd w:id=\"0\"/>' +\n  '</w:p>' +\n  '</w:body></w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\nlastPara.insertOoxml(bookmarkParaOoxml, Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# Fixed #198 Maximum size of 255 char for fields may be a problem for the\n# list of nsURIs to use.\n#\n# The bookmark (\"_GoBack\") that used to sit at the end of the first\n# paragraph is moved to its own (previously empty) trailing paragraph,\n# and a run of the field instruction text that was split across several\n# <w:instrText> runs gets merged into a single run.\n\n$d = $word.ActiveDocument\n\n$nsWord = 'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"'\n\n# 1) Paragraph 1: \"A simple demonstration of a query :\" -- drop the\n#    bookmark that used to live here (it is recreated below, in the\n#    trailing paragraph).\n$firstParaOoxml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document ' + $nsWord + '><w:body>' +\n  '<w:p>' +\n  '<w:r><w:t xml:space=\"preserve\">A simple </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:t>demonstration</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> of a </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:t>query</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  \"<w:r><w:t>$([char]0x00A0):</w:t></w:r>\" +\n  '</w:p>' +\n  '</w:body></w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>'\n$d.Paragraphs.Item(1).Range.InsertXML($firstParaOoxml)\n\n# 2) Paragraph 2: the field-code paragraph. Content/text is unchanged;\n#    only the run split of the instruction text changes: the runs\n#    holding \",\", \"''.sampleTable()\", \",\", \"''.sampleTable()\", \"}\" are\n#    merged into a single run \",''.sampleTable(),''.sampleTable()}\".\n$fieldParaOoxml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document ' + $nsWord + '><w:body>' +\n  '<w:p>' +\n  '<w:r><w:fldChar w:fldCharType=\"begin\"/></w:r>' +\n  '<w:r><w:instrText xml:space=\"preserve\"> </w:instrText></w:r>' +\n  '<w:r><w:instrText>m</w:instrText></w:r>' +\n  '<w:r><w:instrText>:</w:instrText></w:r>' +\n  '<w:r><w:instrText>Sequence</w:instrText></w:r>' +\n  '<w:r><w:instrText>{</w:instrText></w:r>' +\n  \"<w:r><w:instrText>''.</w:instrText></w:r>\" +\n  '<w:r><w:instrText>sample</w:instrText></w:r>' +\n  '<w:r><w:instrText>Table()</w:instrText></w:r>' +\n  \"<w:r><w:instrText>,''.sampleTable(),''.sampleTable()}</w:instrText></w:r>\" +\n  '<w:r><w:instrText xml:space=\"preserve\"> </w:instrText></w:r>' +\n  '<w:r><w:fldChar w:fldCharType=\"end\"/></w:r>' +\n  '</w:p>' +\n  '</w:body></w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>'\n$d.Paragraphs.Item(2).Range.InsertXML($fieldParaOoxml)\n\n# 3) Move the \"_GoBack\" bookmark from the first paragraph to the final\n#    (previously empty) trailing paragraph.\n$d.Bookmarks.Item(\"_GoBack\").Delete()\n$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)\n$d.Bookmarks.Add(\"_GoBack\", $lastPara.Range)\n"}
